# edit.ps1 - applies the "block quotes" paragraph insertion + grammar
# proofErr split + BlockText line-spacing style tweak described by the
# commit diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Locate the paragraph that currently reads:
#      "When you click the Knit button a document will be generated..."
#    and replace it (and its paragraph mark) with four paragraphs:
#      a) the original paragraph (same paraId), now holding the new
#         "You can also include block quotes..." sentence, keeping the
#         <w:lastRenderedPageBreak/> run marker it already had
#      b) a new BlockText-styled paragraph with the long quotation
#      c) a new, empty BodyText paragraph (spacer)
#      d) a new BodyText paragraph with the original "When you click
#         the Knit button..." sentence, except "button" is now split
#         into its own run bracketed by <w:proofErr> gramStart/gramEnd
#         markers (a grammar-check false positive Word flags)
# ---------------------------------------------------------------------

$target = $null
$targetParaId = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "When you click the Knit button a document will be generated*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'When you click the Knit button' paragraph"
}

$replacementXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="0D50E202" w14:textId="77777777" w:rsidR="00957DE8" w:rsidRDefault="00000000">
            <w:pPr>
              <w:pStyle w:val="BodyText"/>
            </w:pPr>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>You can also include block quotes which are hard to demonstrate but they should be in single line spacing:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="BlockText"/>
            </w:pPr>
            <w:r>
              <w:t>Here is lots of text from some quotation that is super important that the reader needs to see in its entirety otherwise they will understand that you&#8217;re a rube and don&#8217;t know what you&#8217;re doing.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="BodyText"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="BodyText"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">When you click the </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>Knit</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>button</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve"> a document will be generated that includes both content as well as the output of any embedded R code chunks within the document. You can embed an R code chunk like this:</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.Range.InsertXML($replacementXml)

# ---------------------------------------------------------------------
# 2) BlockText style: add explicit single-line spacing
#    (w:line="240" w:lineRule="auto") to its paragraph spacing.
# ---------------------------------------------------------------------

$blockTextStyle = $d.Styles("Block Text")
$blockTextStyle.ParagraphFormat.LineSpacingRule = 0
